$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1127.8572
$ws.Range("I43").Value = 1098
$ws.Range("J43").Value = 1146.2307
$ws.Range("K43").Value = 1098
$ws.Range("L43").Value = 1146.2307
$ws.Range("M43").Value = -1029
$ws.Range("N43").Value = -1284.2307
$ws.Range("H51").Value = 2742.4285
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 2639.4
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 2639.4
$ws.Range("M51").Value = -2516
$ws.Range("N51").Value = -3607.4
$ws.Range("H98").Value = 2736862
$ws.Range("I98").Value = 28901.422
$ws.Range("J98").Value = 37037696
$ws.Range("K98").Value = 28901.422
$ws.Range("L98").Value = 37037696
$ws.Range("M98").Value = -27403.422
$ws.Range("N98").Value = -37040692
$ws.Range("H122").Value = 2736862
$ws.Range("I122").Value = 28901.422
$ws.Range("J122").Value = 37037696
$ws.Range("K122").Value = 86704.266
$ws.Range("L122").Value = 111113088
$ws.Range("M122").Value = -84254.266
$ws.Range("N122").Value = -111117988
$ws.Range("H137").Value = 1250.7916
$ws.Range("I137").Value = 1083.125
$ws.Range("J137").Value = 1586.125
$ws.Range("K137").Value = 3249.375
$ws.Range("L137").Value = 4758.375
$ws.Range("M137").Value = -699.375
$ws.Range("N137").Value = -9858.375
$ws.Range("H140").Value = 62387.5
$ws.Range("J140").Value = 88581.25
$ws.Range("L140").Value = 88581.25
$ws.Range("N140").Value = -98941.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4412.778
$ws.Range("I61").Value = 4665.6816
$ws.Range("J61").Value = 3300
$ws.Range("K61").Value = 4665.6816
$ws.Range("L61").Value = 3300
$ws.Range("M61").Value = -4453.6816
$ws.Range("N61").Value = -3724
$ws.Range("H74").Value = 1555
$ws.Range("I74").Value = 806.9231
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 806.9231
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = 67.07690000000002
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 1555
$ws.Range("I77").Value = 806.9231
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 4034.6155
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = 333.3845000000001
$ws.Range("N77").Value = -26236
$ws.Range("H132").Value = 1323.5568
$ws.Range("I132").Value = 1080.4606
$ws.Range("K132").Value = 3241.3818
$ws.Range("M132").Value = -711.3818000000001
$ws.Range("H136").Value = 4412.778
$ws.Range("I136").Value = 4665.6816
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 13997.0448
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -11447.0448
$ws.Range("N136").Value = -15000
$ws.Range("H138").Value = 63400
$ws.Range("J138").Value = 63400
$ws.Range("L138").Value = 63400
$ws.Range("N138").Value = -73680
$ws.Range("H139").Value = 47178.75
$ws.Range("J139").Value = 47178.75
$ws.Range("L139").Value = 47178.75
$ws.Range("N139").Value = -57458.75
$ws.Range("H141").Value = 59786.668
$ws.Range("J141").Value = 62271.43
$ws.Range("L141").Value = 62271.43
$ws.Range("N141").Value = -72631.42999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 62345
$ws.Range("J138").Value = 62345
$ws.Range("L138").Value = 62345
$ws.Range("N138").Value = -72625
$ws.Range("H140").Value = 89750
$ws.Range("J140").Value = 89750
$ws.Range("L140").Value = 89750
$ws.Range("N140").Value = -100110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4766.8
$ws.Range("I31").Value = 5077.3
$ws.Range("K31").Value = 5077.3
$ws.Range("M31").Value = -4782.3
$ws.Range("H34").Value = 4766.8
$ws.Range("I34").Value = 5077.3
$ws.Range("K34").Value = 5077.3
$ws.Range("M34").Value = -4875.3
$ws.Range("H58").Value = 1701.0769
$ws.Range("I58").Value = 1116.8334
$ws.Range("J58").Value = 2201.8572
$ws.Range("K58").Value = 1116.8334
$ws.Range("L58").Value = 2201.8572
$ws.Range("M58").Value = -913.8334
$ws.Range("N58").Value = -2607.8572
$ws.Range("H136").Value = 1701.0769
$ws.Range("I136").Value = 1116.8334
$ws.Range("J136").Value = 2201.8572
$ws.Range("K136").Value = 3350.5002
$ws.Range("L136").Value = 6605.571599999999
$ws.Range("M136").Value = -800.5001999999999
$ws.Range("N136").Value = -11705.5716
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H140").Value = 88369.25
$ws.Range("J140").Value = 88369.25
$ws.Range("L140").Value = 88369.25
$ws.Range("N140").Value = -98729.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 20000440
$ws.Range("I18").Value = 25000344
$ws.Range("J18").Value = 824.2
$ws.Range("K18").Value = 75001032
$ws.Range("L18").Value = 2472.6
$ws.Range("M18").Value = -75000863
$ws.Range("N18").Value = -2810.6
$ws.Range("H131").Value = 6173760.5
$ws.Range("J131").Value = 6945283.5
$ws.Range("L131").Value = 20835850.5
$ws.Range("N131").Value = -20845930.5
$ws.Range("H134").Value = 3315.1082
$ws.Range("I134").Value = 2660.7917
$ws.Range("J134").Value = 4523.077
$ws.Range("K134").Value = 7982.375100000001
$ws.Range("L134").Value = 13569.231
$ws.Range("M134").Value = -2912.375100000001
$ws.Range("N134").Value = -23709.231
$ws.Range("H137").Value = 3042.3333
$ws.Range("I137").Value = 1716
$ws.Range("K137").Value = 5148
$ws.Range("M137").Value = -48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3440815.8
$ws.Range("J11").Value = 5807.6924
$ws.Range("L11").Value = 5807.6924
$ws.Range("N11").Value = -6085.6924
$ws.Range("H21").Value = 1112888.9
$ws.Range("I21").Value = 10000000
$ws.Range("K21").Value = 10000000
$ws.Range("M21").Value = -9999827
$ws.Range("H24").Value = 2000
$ws.Range("J24").Value = 2000
$ws.Range("L24").Value = 2000
$ws.Range("N24").Value = -2346
$ws.Range("H30").Value = 1112888.9
$ws.Range("I30").Value = 10000000
$ws.Range("K30").Value = 10000000
$ws.Range("M30").Value = -9999895
$ws.Range("H138").Value = 68066.664
$ws.Range("J138").Value = 68066.664
$ws.Range("L138").Value = 68066.664
$ws.Range("N138").Value = -78346.664
$ws.Range("H140").Value = 88993
$ws.Range("J140").Value = 88993
$ws.Range("L140").Value = 88993
$ws.Range("N140").Value = -99353

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2968.3076
$ws.Range("I100").Value = 2326.1428
$ws.Range("J100").Value = 3717.5
$ws.Range("K100").Value = 2326.1428
$ws.Range("L100").Value = 3717.5
$ws.Range("M100").Value = -1785.1428
$ws.Range("N100").Value = -4799.5
$ws.Range("H139").Value = 70400
$ws.Range("J139").Value = 70400
$ws.Range("L139").Value = 70400
$ws.Range("N139").Value = -80680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 67260
$ws.Range("J138").Value = 67260
$ws.Range("L138").Value = 67260
$ws.Range("N138").Value = -77540
$ws.Range("H141").Value = 67302.14
$ws.Range("J141").Value = 67302.14
$ws.Range("L141").Value = 67302.14
$ws.Range("N141").Value = -77662.14
